$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: Objetivos -> B/C now hold the "5816812 - João Paulo Alves Silva" text
$ws.Range("B10").Value = "5816812 - João Paulo Alves Silva"
$ws.Range("C10").Value = "5816812 - João Paulo Alves Silva"

# Row 13 gains a label in A ("Programa resumido:") and its B/C text becomes "Semestral"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# Row 14: label becomes "Short syllabus:", B/C cleared
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14:C14").ClearContents()

# Row 15: label becomes "Programa:", B/C become "01/01/2012"
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"
$ws.Rows.Item(15).RowHeight = 120

# Row 16: label becomes "Syllabus:", B/C cleared
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16:C16").ClearContents()

# Row 17: label becomes "Avaliação:" and loses its custom row height
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A17").EntireRow.AutoFit()

# Row 18: label becomes "Método:", B/C become "5816812 - João Paulo Alves Silva"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5816812 - João Paulo Alves Silva"
$ws.Range("C18").Value = "5816812 - João Paulo Alves Silva"
$ws.Rows.Item(18).RowHeight = 60

# Row 19: label becomes "Critério:" (B/C text unchanged)
$ws.Range("A19").Value = "Critério:"

# Row 20: label becomes "Norma de recuperação:" (B/C text unchanged)
$ws.Range("A20").Value = "Norma de recuperação:"

# Row 21: label becomes "Bibliografia:" (B/C text unchanged), height grows
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120

# Row 22: label becomes "Requisitos:", B/C cleared, loses its custom row height
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22:C22").ClearContents()
$ws.Range("A22").EntireRow.AutoFit()

# Row 23: picks up the "LOQ4064..." requirement text that used to live on row 24
$ws.Range("B23").Value = "LOQ4064 -  Engenharia de Processos Quimicos I  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOQ4064 -  Engenharia de Processos Quimicos I  (Requisito fraco)`n"
$ws.Rows.Item(23).RowHeight = 30

# Old row 24 (now redundant) is removed entirely
$ws.Rows.Item(24).Delete()
